$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "28÷9=3, 1"
$t.Cell(1,2).Range.Text = "69÷5=13, 4"
$t.Cell(1,3).Range.Text = "12÷9=1, 3"
$t.Cell(1,4).Range.Text = "84÷4=21, 0"
$t.Cell(1,5).Range.Text = "94÷8=11, 6"
$t.Cell(5,1).Range.Text = "69÷6=11, 3"
$t.Cell(5,2).Range.Text = "38÷4=9, 2"
$t.Cell(5,3).Range.Text = "87÷7=12, 3"
$t.Cell(5,4).Range.Text = "99÷3=33, 0"
$t.Cell(5,5).Range.Text = "86÷7=12, 2"
$t.Cell(9,1).Range.Text = "42÷5=8, 2"
$t.Cell(9,2).Range.Text = "86÷4=21, 2"
$t.Cell(9,3).Range.Text = "16÷8=2, 0"
$t.Cell(9,4).Range.Text = "13÷9=1, 4"
$t.Cell(9,5).Range.Text = "61÷3=20, 1"
$t.Cell(13,1).Range.Text = "80÷8=10, 0"
$t.Cell(13,2).Range.Text = "27÷4=6, 3"
$t.Cell(13,3).Range.Text = "83÷3=27, 2"
$t.Cell(13,4).Range.Text = "54÷3=18, 0"
$t.Cell(13,5).Range.Text = "23÷8=2, 7"
$t.Cell(17,1).Range.Text = "20÷2=10, 0"
$t.Cell(17,2).Range.Text = "76÷8=9, 4"
$t.Cell(17,3).Range.Text = "87÷2=43, 1"
$t.Cell(17,4).Range.Text = "84÷4=21, 0"
$t.Cell(17,5).Range.Text = "96÷6=16, 0"
